$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "ChatGPT" source notes from column D to column C for rows 115-117,
#     and unhide / resize column C to match column B ---
$ws.Range("D115:D117").ClearContents()
$ws.Range("C115:C117").Value = "ChatGPT"

$ws.Columns("C").Hidden = $false
$ws.Columns("C").ColumnWidth = 21

# --- Append new vocabulary rows 118-121 ---
$ws.Range("A118").Value = "Él"
$ws.Range("B118").Value = "Keri"

$ws.Range("A119").Value = "Su"
$ws.Range("B119").Value = "Itsu"

# Note: for row 120 the KAWAIINESE word was entered before the SPANISH term
# so that the shared-string table order matches the source workbook.
$ws.Range("B120").Value = "Pasu"
$ws.Range("A120").Value = "(Gram) Pasado"

$ws.Range("A121").Value = "(Gram) Futuro"
$ws.Range("B121").Value = "Era"

# --- Update the view: scroll down a bit and move the active selection ---
$ws.Range("B122").Select() | Out-Null
